$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.51"
$ws.Range("E2").Value = "'-0.23%"

$ws.Range("D3").Value = "'40.50"
$ws.Range("E3").Value = "'0.90%"

$ws.Range("D4").Value = "'5.013"
$ws.Range("E4").Value = "'-0.20%"

$ws.Range("D5").Value = "'0.07388"
$ws.Range("E5").Value = "'0.12%"

$ws.Range("D6").Value = "'1.578"
$ws.Range("E6").Value = "'1.49%"

$ws.Range("D7").Value = "'0.9240"
$ws.Range("E7").Value = "'0.40%"

$ws.Range("E9").Value = "'-0.01%"

$ws.Range("D10").Value = "'0.1815"
$ws.Range("E10").Value = "'3.48%"

$ws.Range("D11").Value = "'0.04396"
$ws.Range("E11").Value = "'5.13%"

$ws.Range("D12").Value = "'0.08750"

$ws.Range("D14").Value = "'0.001270"
$ws.Range("E14").Value = "'-0.70%"

$ws.Range("D15").Value = "'0.005836"
$ws.Range("E15").Value = "'-0.31%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.342"
$ws.Range("E16").Value = "'-1.20%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.293"
$ws.Range("E17").Value = "'-0.30%"

$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3317"
$ws.Range("E18").Value = "'0.64%"

$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'7.913"
$ws.Range("E19").Value = "'4.25%"

$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1390"
$ws.Range("E20").Value = "'3.40%"

$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2960"
$ws.Range("E21").Value = "'5.38%"

$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.03918"
$ws.Range("E22").Value = "'2.64%"

$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001261"
$ws.Range("E23").Value = "'-1.76%"

$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.003805"
$ws.Range("E24").Value = "'4.16%"

$ws.Range("D25").Value = "'0.0001231"
$ws.Range("E25").Value = "'-4.84%"

$ws.Range("E26").Value = "'-0.24%"

$ws.Range("D38").Value = "'0.02328"
$ws.Range("E38").Value = "'0.98%"

$ws.Range("D39").Value = "'0.05063"
$ws.Range("E39").Value = "'0.94%"

$ws.Range("D40").Value = "'0.005833"
$ws.Range("E40").Value = "'32.07%"

$ws.Range("D41").Value = "'0.007824"
$ws.Range("E41").Value = "'1.39%"

$ws.Range("D42").Value = "'0.1288"
$ws.Range("E42").Value = "'0.94%"

$ws.Range("D43").Value = "'0.007383"
$ws.Range("E43").Value = "'-0.46%"

$ws.Range("D44").Value = "'0.008044"
$ws.Range("E44").Value = "'15.23%"

$ws.Range("D45").Value = "'0.2915"
$ws.Range("E45").Value = "'-8.50%"

$ws.Range("D46").Value = "'0.00006225"
$ws.Range("E46").Value = "'-3.87%"

$ws.Range("E47").Value = "'-0.23%"

$ws.Range("D48").Value = "'0.04795"
$ws.Range("E48").Value = "'-80.96%"

$ws.Range("D49").Value = "'0.004203"
$ws.Range("E49").Value = "'-0.24%"

$ws.Range("E50").Value = "'-0.23%"

$ws.Range("E51").Value = "'-0.23%"
